# Auto-generated edit script for cryptos.xlsx update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.649.01"
$ws.Range("E2").Value = "  +1.05%  "
$ws.Range("D3").Value = "2.656.48"
$ws.Range("E3").Value = "  +2.67%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'594.08"
$ws.Range("E5").Value = "  +1.73%  "
$ws.Range("D6").Value = "'146.69"
$ws.Range("E6").Value = "  -0.37%  "
$ws.Range("E8").Value = "  -0.98%  "
$ws.Range("E9").Value = "  -0.26%  "
$ws.Range("E10").Value = "  -0.40%  "
$ws.Range("E11").Value = "  -0.06%  "
$ws.Range("D12").Value = "'0.355"
$ws.Range("E12").Value = "  +0.39%  "
$ws.Range("D13").Value = "'27.63"
$ws.Range("E13").Value = "  +0.65%  "
$ws.Range("D14").Value = "3.130.65"
$ws.Range("E14").Value = "  +2.63%  "
$ws.Range("D15").Value = "63.424.13"
$ws.Range("E15").Value = "  +0.90%  "
$ws.Range("E16").Value = "  -0.28%  "
$ws.Range("D17").Value = "2.680.49"
$ws.Range("E17").Value = "  +3.00%  "
$ws.Range("D18").Value = "'11.39"
$ws.Range("E18").Value = "  +0.78%  "
$ws.Range("D19").Value = "'342.92"
$ws.Range("E19").Value = "  +0.23%  "
$ws.Range("D20").Value = "'4.36"
$ws.Range("E20").Value = "  -0.75%  "
$ws.Range("D21").Value = "'6.78"
$ws.Range("E21").Value = "  +1.44%  "
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("D23").Value = "'67.96"
$ws.Range("E23").Value = "  +1.12%  "
$ws.Range("B24").Value = "Fetch.AI"
$ws.Range("C24").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D24").Value = "'1.68"
$ws.Range("E24").Value = "  +5.57%  "
$ws.Range("B25").Value = "SuiNetwork"
$ws.Range("C25").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D25").Value = "'1.61"
$ws.Range("E25").Value = "  +12.20%  "
$ws.Range("D26").Value = "'569.89"
$ws.Range("E26").Value = "  +22.83%  "
$ws.Range("E27").Value = "  -0.34%  "
$ws.Range("B28").Value = "InternetComputer(DFINITY)"
$ws.Range("C28").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D28").Value = "'8.52"
$ws.Range("E28").Value = "  +2.46%  "
$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  +0.50%  "
$ws.Range("B30").Value = "Aptos"
$ws.Range("C30").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D30").Value = "'7.97"
$ws.Range("E30").Value = "  +1.08%  "
$ws.Range("E31").Value = "  +4.04%  "
$ws.Range("D32").Value = "'1.79"
$ws.Range("E32").Value = "  +11.53%  "
$ws.Range("D33").Value = "0.0₃0816"
$ws.Range("E33").Value = "  -0.72%  "
$ws.Range("D34").Value = "'175.25"
$ws.Range("E34").Value = "  +0.18%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").Value = "'4.81"
$ws.Range("E36").Value = "  +6.03%  "
$ws.Range("B37").Value = "PolygonEcosystemToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D37").Value = "'0.401"
$ws.Range("E37").Value = "  -0.60%  "
$ws.Range("D38").Value = "'19.15"
$ws.Range("E38").Value = "  +0.36%  "
$ws.Range("D39").Value = "'1.75"
$ws.Range("E39").Value = "  +2.80%  "
$ws.Range("D40").Value = "'170.60"
$ws.Range("E40").Value = "  +7.26%  "
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("D42").Value = "'40.51"
$ws.Range("E42").Value = "  +2.91%  "
$ws.Range("D43").Value = "'3.76"
$ws.Range("E43").Value = "  -0.21%  "
$ws.Range("D44").Value = "'21.72"
$ws.Range("E44").Value = "  +2.36%  "
$ws.Range("D45").Value = "'0.630"
$ws.Range("E45").Value = "  -1.32%  "
$ws.Range("D46").Value = "'0.0556"
$ws.Range("E46").Value = "  +2.61%  "
$ws.Range("E47").Value = "  +1.52%  "
$ws.Range("D48").Value = "'0.0960"
$ws.Range("E48").Value = "  -0.81%  "
$ws.Range("E49").Value = "  +1.39%  "
$ws.Range("D50").Value = "'1.73"
$ws.Range("E50").Value = "  +0.82%  "
$ws.Range("E51").Value = "  +12.59%  "

# Reset quote-prefix styling introduced by forcing text above, so cell
# formatting matches the original (unstyled) state.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").Style = "Normal"
